$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.699.27"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.851.35"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'312.63"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.4265"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").Value = "'0.3640"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").Value = "'0.07288"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'20.58"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "1.855.79"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "'5.313"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'6.508"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'0.06883"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'79.76"
$ws.Range("E18").Value = "  +3.55%  "
$ws.Range("D19").Value = "'0.000009034"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "27.691.39"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'4.965"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'10.38"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("D25").Value = "2.107.65"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").Value = "'1.965"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").Value = "'153.51"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "'18.81"
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").Value = "'122.29"
$ws.Range("E29").Value = "  +10.71%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +10.90%  "
$ws.Range("D32").Value = "'0.08892"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'0.7608"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").Value = "'2.968"
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").Value = "'4.517"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'0.05369"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'0.01928"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "'2.813"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("D41").Value = "'0.5049"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'0.1645"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'6.762"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'8.341"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").Value = "'0.06542"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").Value = "'10.29"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "'104.83"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "'0.4648"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'1.620"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "'64.36"
$ws.Range("E51").Value = "  +0.00%  "
